$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The crawler now records only the single most-recent date (and drops the
# per-tenor yield columns/rows that used to be scraped alongside it), so
# everything except A1 goes away and A1 becomes "9-9-2023".

# 1) Wipe out the old header row (B1:L1) and the extra date rows (2-4),
#    along with column B1:L4, leaving just A1.
$ws.Range("A1:L4").ClearContents()

# 2) Write the new value into A1. Format it as text first so Excel stores
#    the literal string "9-9-2023" instead of auto-converting it into a
#    date serial number (it still reads/parses as text, same as before).
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "9-9-2023"

# 3) Drop the temporary text format again so the cell keeps the workbook's
#    default (General) style, matching the rest of the (now-empty) sheet.
$ws.Range("A1").ClearFormats()
